$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2098360655737705
$ws.Range("C2").Value = 0.5245901639344263
$ws.Range("J2").Value = 0.03606557377049181
$ws.Range("P2").Value = 0.1213114754098361
$ws.Range("S2").Value = 0.1081967213114754
$ws.Range("B3").Value = 0.01226993865030675
$ws.Range("C3").Value = 0.01226993865030675
$ws.Range("J3").Value = 0.0245398773006135
$ws.Range("P3").Value = 0.754601226993865
$ws.Range("S3").Value = 0.196319018404908
$ws.Range("P4").Value = 0.6521739130434783
$ws.Range("S4").Value = 0.3478260869565217
$ws.Range("B6").Value = 0.07627118644067797
$ws.Range("D6").Value = 0.008474576271186441
$ws.Range("F6").Value = 0.1016949152542373
$ws.Range("J6").Value = 0.1991525423728814
$ws.Range("O6").Value = 0.03813559322033899
$ws.Range("Q6").Value = 0.1313559322033898
$ws.Range("R6").Value = 0.0847457627118644
$ws.Range("S6").Value = 0.3601694915254237
$ws.Range("B7").Value = 0.08854166666666667
$ws.Range("D7").Value = 0.02604166666666667
$ws.Range("F7").Value = 0.05208333333333334
$ws.Range("J7").Value = 0.125
$ws.Range("O7").Value = 0.02083333333333333
$ws.Range("Q7").Value = 0.1614583333333333
$ws.Range("R7").Value = 0.078125
$ws.Range("S7").Value = 0.4479166666666667
$ws.Range("B8").Value = 0.1067538126361656
$ws.Range("D8").Value = 0.02178649237472767
$ws.Range("F8").Value = 0.07625272331154684
$ws.Range("J8").Value = 0.1089324618736384
$ws.Range("O8").Value = 0.01742919389978214
$ws.Range("Q8").Value = 0.1459694989106754
$ws.Range("R8").Value = 0.08278867102396514
$ws.Range("S8").Value = 0.4400871459694989
$ws.Range("B9").Value = 0.1382488479262673
$ws.Range("D9").Value = 0.02304147465437788
$ws.Range("E9").Value = 0.004608294930875576
$ws.Range("F9").Value = 0.05990783410138249
$ws.Range("J9").Value = 0.07834101382488479
$ws.Range("O9").Value = 0.0184331797235023
$ws.Range("Q9").Value = 0.1658986175115207
$ws.Range("R9").Value = 0.1105990783410138
$ws.Range("S9").Value = 0.4009216589861751
$ws.Range("B10").Value = 0.1060869565217391
$ws.Range("D10").Value = 0.02260869565217391
$ws.Range("E10").Value = 0.0008695652173913044
$ws.Range("F10").Value = 0.07739130434782608
$ws.Range("J10").Value = 0.09913043478260869
$ws.Range("O10").Value = 0.01391304347826087
$ws.Range("Q10").Value = 0.2095652173913043
$ws.Range("R10").Value = 0.09130434782608696
$ws.Range("S10").Value = 0.3791304347826087
$ws.Range("G11").Value = 0.1169811320754717
$ws.Range("J11").Value = 0.06792452830188679
$ws.Range("K11").Value = 0.1584905660377358
$ws.Range("L11").Value = 0.6528301886792452
$ws.Range("S11").Value = 0.003773584905660377
$ws.Range("G12").Value = 0.7611111111111111
$ws.Range("J12").Value = 0.1888888888888889
$ws.Range("L12").Value = 0.02777777777777778
$ws.Range("S12").Value = 0.02222222222222222
$ws.Range("F15").Value = 0.02439024390243903
$ws.Range("H15").Value = 0.1414634146341463
$ws.Range("I15").Value = 0.07804878048780488
$ws.Range("J15").Value = 0.3853658536585366
$ws.Range("K15").Value = 0.04878048780487805
$ws.Range("M15").Value = 0.004878048780487805
$ws.Range("N15").Value = 0.004878048780487805
$ws.Range("O15").Value = 0.07317073170731707
$ws.Range("S15").Value = 0.2390243902439024
$ws.Range("F16").Value = 0.01092896174863388
$ws.Range("H16").Value = 0.2295081967213115
$ws.Range("I16").Value = 0.07650273224043716
$ws.Range("J16").Value = 0.3497267759562842
$ws.Range("K16").Value = 0.1202185792349727
$ws.Range("M16").Value = 0.01092896174863388
$ws.Range("O16").Value = 0.0546448087431694
$ws.Range("S16").Value = 0.1475409836065574
$ws.Range("F17").Value = 0.01485148514851485
$ws.Range("H17").Value = 0.1806930693069307
$ws.Range("I17").Value = 0.1014851485148515
$ws.Range("J17").Value = 0.3886138613861386
$ws.Range("K17").Value = 0.1163366336633663
$ws.Range("M17").Value = 0.02475247524752475
$ws.Range("O17").Value = 0.05693069306930693
$ws.Range("S17").Value = 0.1163366336633663
$ws.Range("F18").Value = 0.02475247524752475
$ws.Range("H18").Value = 0.1584158415841584
$ws.Range("I18").Value = 0.0891089108910891
$ws.Range("J18").Value = 0.4554455445544555
$ws.Range("K18").Value = 0.07920792079207921
$ws.Range("M18").Value = 0.0198019801980198
$ws.Range("N18").Value = 0.004950495049504951
$ws.Range("O18").Value = 0.0594059405940594
$ws.Range("S18").Value = 0.1089108910891089
$ws.Range("F19").Value = 0.01115537848605578
$ws.Range("H19").Value = 0.2294820717131474
$ws.Range("I19").Value = 0.1035856573705179
$ws.Range("J19").Value = 0.349003984063745
$ws.Range("K19").Value = 0.101195219123506
$ws.Range("M19").Value = 0.02310756972111554
$ws.Range("O19").Value = 0.06852589641434263
$ws.Range("S19").Value = 0.1131474103585657

Write-Host "Applied 105 cell updates"
